$p = $ppt.ActivePresentation

# New slide 3 uses the same "Title Text Plaid Footers" layout as slide 1
$layout = $p.Slides.Item(1).CustomLayout
$s = $p.Slides.AddSlide(3, $layout)

# Title
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Test"

# Turn on the slide-number placeholder (layout defines one at idx=10)
$s.HeadersFooters.SlideNumber.Visible = $true

# Re-order so the slide matches Title / Slide Number / Content
$s.Shapes.Item(3).ZOrder(3)

# Match shape names used by the authored slide
$s.Shapes.Item(2).Name = "Slide Number Placeholder 2"
$s.Shapes.Item(3).Name = "Content Placeholder 3"
